$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append at the end of the data (matching the existing
# date/nuovi pos./somma mobile 7gg./somma mobile 7gg. per 100mila abitanti
# series), continuing the update through 13/05 (serial 44329).
$newRows = @(
    @(252, 44326, 0, 2, 33.27233405423391),
    @(253, 44327, 0, 1, 16.63616702711695),
    @(254, 44328, 0, 1, 16.63616702711695),
    @(255, 44329, 0, 1, 16.63616702711695)
)

# Use the last existing data row as the style template for column A
# (date formatted, bordered, centered style already used throughout
# the column) and copy it down, then overwrite the values.
$templateCell = $ws.Cells.Item(251, 1)

foreach ($row in $newRows) {
    $r = $row[0]

    $dateCell = $ws.Cells.Item($r, 1)
    $templateCell.Copy($dateCell)
    $dateCell.Value2 = $row[1]

    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}
